$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 44253
$ws.Range("M2").Value = 12
$ws.Range("O2").Value = 200000
$ws.Range("P2").Value = 190000
$ws.Range("S2").Value = 190000

# Row 4 updates
$ws.Range("D4").Value = 44672
$ws.Range("M4").Value = 8
$ws.Range("O4").Value = 180000
$ws.Range("P4").Value = 180000
$ws.Range("S4").Value = 180000
